# Real-Time forecast dataset update: shift timestamps to 2024-06-19 and
# refresh Power/Energy columns (B-F) for rows 24-63 with real-time production data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps in column A forward by 8 days (2024-06-11 -> 2024-06-19)
$ws.Range("A2").Value = "2024-06-19 00:00:00+00:00"
$ws.Range("A3").Value = "2024-06-19 00:15:00+00:00"
$ws.Range("A4").Value = "2024-06-19 00:30:00+00:00"
$ws.Range("A5").Value = "2024-06-19 00:45:00+00:00"
$ws.Range("A6").Value = "2024-06-19 01:00:00+00:00"
$ws.Range("A7").Value = "2024-06-19 01:15:00+00:00"
$ws.Range("A8").Value = "2024-06-19 01:30:00+00:00"
$ws.Range("A9").Value = "2024-06-19 01:45:00+00:00"
$ws.Range("A10").Value = "2024-06-19 02:00:00+00:00"
$ws.Range("A11").Value = "2024-06-19 02:15:00+00:00"
$ws.Range("A12").Value = "2024-06-19 02:30:00+00:00"
$ws.Range("A13").Value = "2024-06-19 02:45:00+00:00"
$ws.Range("A14").Value = "2024-06-19 03:00:00+00:00"
$ws.Range("A15").Value = "2024-06-19 03:15:00+00:00"
$ws.Range("A16").Value = "2024-06-19 03:30:00+00:00"
$ws.Range("A17").Value = "2024-06-19 03:45:00+00:00"
$ws.Range("A18").Value = "2024-06-19 04:00:00+00:00"
$ws.Range("A19").Value = "2024-06-19 04:15:00+00:00"
$ws.Range("A20").Value = "2024-06-19 04:30:00+00:00"
$ws.Range("A21").Value = "2024-06-19 04:45:00+00:00"
$ws.Range("A22").Value = "2024-06-19 05:00:00+00:00"
$ws.Range("A23").Value = "2024-06-19 05:15:00+00:00"
$ws.Range("A24").Value = "2024-06-19 05:30:00+00:00"
$ws.Range("A25").Value = "2024-06-19 05:45:00+00:00"
$ws.Range("A26").Value = "2024-06-19 06:00:00+00:00"
$ws.Range("A27").Value = "2024-06-19 06:15:00+00:00"
$ws.Range("A28").Value = "2024-06-19 06:30:00+00:00"
$ws.Range("A29").Value = "2024-06-19 06:45:00+00:00"
$ws.Range("A30").Value = "2024-06-19 07:00:00+00:00"
$ws.Range("A31").Value = "2024-06-19 07:15:00+00:00"
$ws.Range("A32").Value = "2024-06-19 07:30:00+00:00"
$ws.Range("A33").Value = "2024-06-19 07:45:00+00:00"
$ws.Range("A34").Value = "2024-06-19 08:00:00+00:00"
$ws.Range("A35").Value = "2024-06-19 08:15:00+00:00"
$ws.Range("A36").Value = "2024-06-19 08:30:00+00:00"
$ws.Range("A37").Value = "2024-06-19 08:45:00+00:00"
$ws.Range("A38").Value = "2024-06-19 09:00:00+00:00"
$ws.Range("A39").Value = "2024-06-19 09:15:00+00:00"
$ws.Range("A40").Value = "2024-06-19 09:30:00+00:00"
$ws.Range("A41").Value = "2024-06-19 09:45:00+00:00"
$ws.Range("A42").Value = "2024-06-19 10:00:00+00:00"
$ws.Range("A43").Value = "2024-06-19 10:15:00+00:00"
$ws.Range("A44").Value = "2024-06-19 10:30:00+00:00"
$ws.Range("A45").Value = "2024-06-19 10:45:00+00:00"
$ws.Range("A46").Value = "2024-06-19 11:00:00+00:00"
$ws.Range("A47").Value = "2024-06-19 11:15:00+00:00"
$ws.Range("A48").Value = "2024-06-19 11:30:00+00:00"
$ws.Range("A49").Value = "2024-06-19 11:45:00+00:00"
$ws.Range("A50").Value = "2024-06-19 12:00:00+00:00"
$ws.Range("A51").Value = "2024-06-19 12:15:00+00:00"
$ws.Range("A52").Value = "2024-06-19 12:30:00+00:00"
$ws.Range("A53").Value = "2024-06-19 12:45:00+00:00"
$ws.Range("A54").Value = "2024-06-19 13:00:00+00:00"
$ws.Range("A55").Value = "2024-06-19 13:15:00+00:00"
$ws.Range("A56").Value = "2024-06-19 13:30:00+00:00"
$ws.Range("A57").Value = "2024-06-19 13:45:00+00:00"
$ws.Range("A58").Value = "2024-06-19 14:00:00+00:00"
$ws.Range("A59").Value = "2024-06-19 14:15:00+00:00"
$ws.Range("A60").Value = "2024-06-19 14:30:00+00:00"
$ws.Range("A61").Value = "2024-06-19 14:45:00+00:00"
$ws.Range("A62").Value = "2024-06-19 15:00:00+00:00"
$ws.Range("A63").Value = "2024-06-19 15:15:00+00:00"
$ws.Range("A64").Value = "2024-06-19 15:30:00+00:00"
$ws.Range("A65").Value = "2024-06-19 15:45:00+00:00"
$ws.Range("A66").Value = "2024-06-19 16:00:00+00:00"
$ws.Range("A67").Value = "2024-06-19 16:15:00+00:00"
$ws.Range("A68").Value = "2024-06-19 16:30:00+00:00"
$ws.Range("A69").Value = "2024-06-19 16:45:00+00:00"
$ws.Range("A70").Value = "2024-06-19 17:00:00+00:00"
$ws.Range("A71").Value = "2024-06-19 17:15:00+00:00"
$ws.Range("A72").Value = "2024-06-19 17:30:00+00:00"
$ws.Range("A73").Value = "2024-06-19 17:45:00+00:00"
$ws.Range("A74").Value = "2024-06-19 18:00:00+00:00"
$ws.Range("A75").Value = "2024-06-19 18:15:00+00:00"
$ws.Range("A76").Value = "2024-06-19 18:30:00+00:00"
$ws.Range("A77").Value = "2024-06-19 18:45:00+00:00"
$ws.Range("A78").Value = "2024-06-19 19:00:00+00:00"
$ws.Range("A79").Value = "2024-06-19 19:15:00+00:00"
$ws.Range("A80").Value = "2024-06-19 19:30:00+00:00"
$ws.Range("A81").Value = "2024-06-19 19:45:00+00:00"
$ws.Range("A82").Value = "2024-06-19 20:00:00+00:00"
$ws.Range("A83").Value = "2024-06-19 20:15:00+00:00"
$ws.Range("A84").Value = "2024-06-19 20:30:00+00:00"
$ws.Range("A85").Value = "2024-06-19 20:45:00+00:00"
$ws.Range("A86").Value = "2024-06-19 21:00:00+00:00"
$ws.Range("A87").Value = "2024-06-19 21:15:00+00:00"
$ws.Range("A88").Value = "2024-06-19 21:30:00+00:00"
$ws.Range("A89").Value = "2024-06-19 21:45:00+00:00"
$ws.Range("A90").Value = "2024-06-19 22:00:00+00:00"
$ws.Range("A91").Value = "2024-06-19 22:15:00+00:00"
$ws.Range("A92").Value = "2024-06-19 22:30:00+00:00"
$ws.Range("A93").Value = "2024-06-19 22:45:00+00:00"
$ws.Range("A94").Value = "2024-06-19 23:00:00+00:00"
$ws.Range("A95").Value = "2024-06-19 23:15:00+00:00"
$ws.Range("A96").Value = "2024-06-19 23:30:00+00:00"
$ws.Range("A97").Value = "2024-06-19 23:45:00+00:00"

# Update forecast values (columns B-F) for rows 24-63 with real-time production data
$ws.Range("B24").Value = 689.099866231283
$ws.Range("C24").Value = 0.0006890998662312
$ws.Range("E24").Value = 0.0003445499331156
$ws.Range("F24").Value = 0.00008613748327891039
$ws.Range("B25").Value = 7859.702987670898
$ws.Range("C25").Value = 0.0078597029876708
$ws.Range("D25").Value = 0.0006890998662312
$ws.Range("E25").Value = 0.004274401426951
$ws.Range("F25").Value = 0.0010686003567377
$ws.Range("B26").Value = 21855.46544392904
$ws.Range("C26").Value = 0.021855465443929
$ws.Range("D26").Value = 0.0078597029876708
$ws.Range("E26").Value = 0.0148575842157999
$ws.Range("F26").Value = 0.0037143960539499
$ws.Range("B27").Value = 47716.76338704427
$ws.Range("C27").Value = 0.0477167633870442
$ws.Range("D27").Value = 0.021855465443929
$ws.Range("E27").Value = 0.0347861144154866
$ws.Range("F27").Value = 0.0086965286038716
$ws.Range("B28").Value = 66857.48278808594
$ws.Range("C28").Value = 0.0668574827880859
$ws.Range("D28").Value = 0.0477167633870442
$ws.Range("E28").Value = 0.057287123087565
$ws.Range("F28").Value = 0.0143217807718912
$ws.Range("B29").Value = 86289.01322428384
$ws.Range("C29").Value = 0.0862890132242838
$ws.Range("D29").Value = 0.0668574827880859
$ws.Range("E29").Value = 0.0765732480061848
$ws.Range("F29").Value = 0.0191433120015461
$ws.Range("B30").Value = 130161.9264322917
$ws.Range("C30").Value = 0.1301619264322915
$ws.Range("D30").Value = 0.0862890132242838
$ws.Range("E30").Value = 0.1082254698282877
$ws.Range("F30").Value = 0.0270563674570719
$ws.Range("B31").Value = 229625.5572102864
$ws.Range("C31").Value = 0.2296255572102864
$ws.Range("D31").Value = 0.1301619264322915
$ws.Range("E31").Value = 0.179893741821289
$ws.Range("F31").Value = 0.0449734354553222
$ws.Range("B32").Value = 386699.9888509115
$ws.Range("C32").Value = 0.3866999888509115
$ws.Range("D32").Value = 0.2296255572102864
$ws.Range("E32").Value = 0.3081627730305988
$ws.Range("F32").Value = 0.0770406932576496
$ws.Range("B33").Value = 594515.7807617188
$ws.Range("C33").Value = 0.5945157807617187
$ws.Range("D33").Value = 0.3866999888509115
$ws.Range("E33").Value = 0.4906078848063151
$ws.Range("F33").Value = 0.1226519712015787
$ws.Range("B34").Value = 803619.599609375
$ws.Range("C34").Value = 0.8036195996093749
$ws.Range("D34").Value = 0.5945157807617187
$ws.Range("E34").Value = 0.6990676901855467
$ws.Range("F34").Value = 0.1747669225463866
$ws.Range("B35").Value = 1028912.063802083
$ws.Range("C35").Value = 1.028912063802083
$ws.Range("D35").Value = 0.8036195996093749
$ws.Range("E35").Value = 0.9162658317057291
$ws.Range("F35").Value = 0.2290664579264322
$ws.Range("B36").Value = 1262598.486653646
$ws.Range("C36").Value = 1.262598486653646
$ws.Range("D36").Value = 1.028912063802083
$ws.Range("E36").Value = 1.145755275227865
$ws.Range("F36").Value = 0.2864388188069661
$ws.Range("B37").Value = 1464948.033203125
$ws.Range("C37").Value = 1.464948033203125
$ws.Range("D37").Value = 1.262598486653646
$ws.Range("E37").Value = 1.363773259928386
$ws.Range("F37").Value = 0.3409433149820963
$ws.Range("B38").Value = 1688063.014322917
$ws.Range("C38").Value = 1.688063014322917
$ws.Range("D38").Value = 1.464948033203125
$ws.Range("E38").Value = 1.576505523763021
$ws.Range("F38").Value = 0.3941263809407551
$ws.Range("B39").Value = 1903729.349609375
$ws.Range("C39").Value = 1.903729349609375
$ws.Range("D39").Value = 1.688063014322917
$ws.Range("E39").Value = 1.795896181966146
$ws.Range("F39").Value = 0.4489740454915365
$ws.Range("B40").Value = 2111411.909505208
$ws.Range("C40").Value = 2.111411909505208
$ws.Range("D40").Value = 1.903729349609375
$ws.Range("E40").Value = 2.007570629557291
$ws.Range("F40").Value = 0.5018926573893228
$ws.Range("B41").Value = 2301467.50390625
$ws.Range("C41").Value = 2.30146750390625
$ws.Range("D41").Value = 2.111411909505208
$ws.Range("E41").Value = 2.206439706705729
$ws.Range("F41").Value = 0.5516099266764323
$ws.Range("B42").Value = 2437485.73046875
$ws.Range("C42").Value = 2.43748573046875
$ws.Range("D42").Value = 2.30146750390625
$ws.Range("E42").Value = 2.3694766171875
$ws.Range("F42").Value = 0.592369154296875
$ws.Range("B43").Value = 2569529.276692708
$ws.Range("C43").Value = 2.569529276692708
$ws.Range("D43").Value = 2.43748573046875
$ws.Range("E43").Value = 2.503507503580729
$ws.Range("F43").Value = 0.6258768758951823
$ws.Range("B44").Value = 2716929.630208333
$ws.Range("C44").Value = 2.716929630208333
$ws.Range("D44").Value = 2.569529276692708
$ws.Range("E44").Value = 2.643229453450521
$ws.Range("F44").Value = 0.6608073633626301
$ws.Range("B45").Value = 2877964.766927083
$ws.Range("C45").Value = 2.877964766927083
$ws.Range("D45").Value = 2.716929630208333
$ws.Range("E45").Value = 2.797447198567709
$ws.Range("F45").Value = 0.699361799641927
$ws.Range("B46").Value = 2954226.1484375
$ws.Range("C46").Value = 2.9542261484375
$ws.Range("D46").Value = 2.877964766927083
$ws.Range("E46").Value = 2.916095457682292
$ws.Range("F46").Value = 0.7290238644205729
$ws.Range("B47").Value = 2937693.569010417
$ws.Range("C47").Value = 2.937693569010416
$ws.Range("D47").Value = 2.9542261484375
$ws.Range("E47").Value = 2.945959858723958
$ws.Range("F47").Value = 0.7364899646809895
$ws.Range("B48").Value = 3020503.787760417
$ws.Range("C48").Value = 3.020503787760417
$ws.Range("D48").Value = 2.937693569010416
$ws.Range("E48").Value = 2.979098678385417
$ws.Range("F48").Value = 0.7447746695963542
$ws.Range("B49").Value = 3196022.669270833
$ws.Range("C49").Value = 3.196022669270834
$ws.Range("D49").Value = 3.020503787760417
$ws.Range("E49").Value = 3.108263228515625
$ws.Range("F49").Value = 0.7770658071289062
$ws.Range("B50").Value = 3195557.708333333
$ws.Range("C50").Value = 3.195557708333333
$ws.Range("D50").Value = 3.196022669270834
$ws.Range("E50").Value = 3.195790188802083
$ws.Range("F50").Value = 0.7989475472005209
$ws.Range("B51").Value = 3245345.85546875
$ws.Range("C51").Value = 3.24534585546875
$ws.Range("D51").Value = 3.195557708333333
$ws.Range("E51").Value = 3.220451781901041
$ws.Range("F51").Value = 0.8051129454752602
$ws.Range("B52").Value = 3225369.536458333
$ws.Range("C52").Value = 3.225369536458333
$ws.Range("D52").Value = 3.24534585546875
$ws.Range("E52").Value = 3.235357695963542
$ws.Range("F52").Value = 0.8088394239908854
$ws.Range("B53").Value = 3313122.616536458
$ws.Range("C53").Value = 3.313122616536458
$ws.Range("D53").Value = 3.225369536458333
$ws.Range("E53").Value = 3.269246076497396
$ws.Range("F53").Value = 0.817311519124349
$ws.Range("B54").Value = 0.0
$ws.Range("C54").Value = 0.0
$ws.Range("D54").Value = 3.313122616536458
$ws.Range("E54").Value = 1.656561308268229
$ws.Range("F54").Value = 0.4141403270670572
$ws.Range("B55").Value = 0.0
$ws.Range("C55").Value = 0.0
$ws.Range("D55").Value = 0.0
$ws.Range("E55").Value = 0.0
$ws.Range("F55").Value = 0.0
$ws.Range("B56").Value = 0.0
$ws.Range("C56").Value = 0.0
$ws.Range("D56").Value = 0.0
$ws.Range("E56").Value = 0.0
$ws.Range("F56").Value = 0.0
$ws.Range("B57").Value = 0.0
$ws.Range("C57").Value = 0.0
$ws.Range("D57").Value = 0.0
$ws.Range("E57").Value = 0.0
$ws.Range("F57").Value = 0.0
$ws.Range("B58").Value = 0.0
$ws.Range("C58").Value = 0.0
$ws.Range("D58").Value = 0.0
$ws.Range("E58").Value = 0.0
$ws.Range("F58").Value = 0.0
$ws.Range("B59").Value = 0.0
$ws.Range("C59").Value = 0.0
$ws.Range("D59").Value = 0.0
$ws.Range("E59").Value = 0.0
$ws.Range("F59").Value = 0.0
$ws.Range("B60").Value = 0.0
$ws.Range("C60").Value = 0.0
$ws.Range("D60").Value = 0.0
$ws.Range("E60").Value = 0.0
$ws.Range("F60").Value = 0.0
$ws.Range("B61").Value = 0.0
$ws.Range("C61").Value = 0.0
$ws.Range("D61").Value = 0.0
$ws.Range("E61").Value = 0.0
$ws.Range("F61").Value = 0.0
$ws.Range("B62").Value = 0.0
$ws.Range("C62").Value = 0.0
$ws.Range("D62").Value = 0.0
$ws.Range("E62").Value = 0.0
$ws.Range("F62").Value = 0.0
$ws.Range("D63").Value = 0.0
$ws.Range("E63").Value = 0.0
$ws.Range("F63").Value = 0.0
